$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple "want to go" count (column F) bumps on existing rows ---
    $ws.Cells.Item(7, 6).Value  = 1698
    $ws.Cells.Item(8, 6).Value  = 26
    $ws.Cells.Item(9, 6).Value  = 692
    $ws.Cells.Item(13, 6).Value = 102
    $ws.Cells.Item(14, 6).Value = 228
    $ws.Cells.Item(19, 6).Value = 3930
    $ws.Cells.Item(22, 6).Value = 443
    $ws.Cells.Item(23, 6).Value = 373
    $ws.Cells.Item(24, 6).Value = 831
    $ws.Cells.Item(25, 6).Value = 571
    $ws.Cells.Item(28, 6).Value = 1727
    $ws.Cells.Item(29, 6).Value = 22
    $ws.Cells.Item(30, 6).Value = 28

    # --- Insert a new event row at row 31, pushing the old 31/32 down to 32/33 ---
    $ws.Rows.Item(31).Insert()

    # Copy formatting from the row that landed on 32 (the old row 31) onto the
    # freshly inserted row so column A keeps the bold/bordered "index" style.
    $ws.Cells.Item(32, 1).Copy($ws.Cells.Item(31, 1))

    $ws.Cells.Item(31, 1).Value = 30

    # Plain "YYYY-MM-DD" text gets auto-parsed into a date serial by a bare
    # Value assignment; force it to stay literal text (matching the sheet's
    # existing inlineStr date cells) and then drop the now-unneeded format.
    $ws.Cells.Item(31, 2).NumberFormat = "@"
    $ws.Cells.Item(31, 2).Value = "2024-05-03"
    $ws.Cells.Item(31, 2).ClearFormats()

    $ws.Cells.Item(31, 3).Value = "江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会"
    $ws.Cells.Item(31, 4).Value = "前湖大道欣悦湖体育馆 欣悦湖体育馆"
    $ws.Cells.Item(31, 5).Value = "2024.05.03 09:30-05.03 17:30"
    $ws.Cells.Item(31, 6).Value = 29
    $ws.Cells.Item(31, 7).Value = 188
    $ws.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83497"
    $ws.Cells.Item(31, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg"

    # The event that used to be row 31 is now row 32; its "want to go" count bumped too.
    $ws.Cells.Item(32, 6).Value = 175
}
